# Adicionado data de controle
# Adds two new tracking columns (DataCriacao / DataModificacao) and a new
# task row ("Mestre do Capitalismo") to the kanban worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:G1 -------------------------------------------------
$ws.Range("F1").Value = "DataCriacao"
$ws.Range("G1").Value = "DataModificacao"

# Copy the header formatting (bold, border, centered/top aligned) from the
# existing header row onto the two new header cells.
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Touch F2:G22 so they exist as blank cells (matches the row-by-row
# export of the source tool, which writes an empty cell for every column in
# every existing row once new columns are introduced) ----------------------
$ws.Range("F2:G22").Style = "Normal"

# --- New row 23: "Mestre do Capitalismo" ------------------------------------
$ws.Range("A23").Value = "Mestre do Capitalismo"
$ws.Range("B23").Value = "A FAZER"
$ws.Range("C23").Value = "Finanças"
$ws.Range("D23").Value = "Alta"
$ws.Range("E23").Value = "copiar videos:False"
$ws.Range("F23").Value = "2025-04-16 21:43:05"
$ws.Range("G23").Value = "2025-04-16 21:43:05"
